$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: ".Value" getter is unreliable in this runtime (returns a
# reflection-looking placeholder string instead of the real value), so all
# reads below use ".Value2" instead. ".Value" as a *setter* works fine.

# --- Read current (pre-edit) values for the cells that need to swap ---
$a7 = $ws.Range("A7").Value2
$a9 = $ws.Range("A9").Value2

$i7 = $ws.Range("I7").Value2
$i9 = $ws.Range("I9").Value2

$q7 = $ws.Range("Q7").Value2
$q9 = $ws.Range("Q9").Value2

$r7 = $ws.Range("R7").Value2
$r9 = $ws.Range("R9").Value2

$ac7 = $ws.Range("AC7").Value2
$ac9 = $ws.Range("AC9").Value2

# --- A (numeric Id) : simple swap ---
$ws.Range("A7").Value = $a9
$ws.Range("A9").Value = $a7

# --- I (Antal) : stored as text in the original file even though the
# content looks numeric ("20"/"15"). Force text storage by setting the
# NumberFormat to "@" before assigning, then reset the display style back
# to Normal so no stray formatting is left behind on the cell. ---
$ws.Range("I7").NumberFormat = "@"
$ws.Range("I7").Value = [string]$i9
$ws.Range("I7").Style = "Normal"

$ws.Range("I9").NumberFormat = "@"
$ws.Range("I9").Value = [string]$i7
$ws.Range("I9").Style = "Normal"

# --- Q / R (coordinates) : simple numeric swap ---
$ws.Range("Q7").Value = $q9
$ws.Range("Q9").Value = $q7

$ws.Range("R7").Value = $r9
$ws.Range("R9").Value = $r7

# --- AC (public comment) : simple text swap ---
$ws.Range("AC7").Value = $ac9
$ws.Range("AC9").Value = $ac7
